$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") — reuse the header style from H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-7: column I is always 1, column J mirrors column H
$hValues = @(2, 5, 4, 4, 3, 2)
for ($i = 0; $i -lt 6; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = 1
    $ws.Cells.Item($row, 10).Value = $hValues[$i]
}
